# Cut down the number of routes to get davinci-003 support.
#
# Removes the "GET Available Flights", "GET Hotel Availability",
# "POST Create New Travel Package", "GET All Travel Packages" sections
# in full, and the heading/Route/Description (but not the Action bullet)
# of the "GET Travel Package Details" section.

$d = $word.ActiveDocument

$startPara = $d.Paragraphs.Item(22)   # "GET Method - Get Available Flights"
$endPara   = $d.Paragraphs.Item(40)   # "Description: Retrieves detailed information about a specific travel package."

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
